# Properly calculate tax and discount in quotation
#
# The "taxRateId" column (column E) is no longer imported directly — tax is
# now calculated elsewhere — so remove it entirely from the product import
# template. Deleting the whole column shifts quantity/costPrice/.../unitId
# one column to the left (E..N -> E..M) and leaves the selection on the
# former taxRateId header cell (now "quantity", E1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()
$ws.Range("E1").Select()
